$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 8124.3887
$ws.Range("I62").Value = 8514.058999999999
$ws.Range("J62").Value = 1500
$ws.Range("K62").Value = 8514.058999999999
$ws.Range("L62").Value = 1500
$ws.Range("M62").Value = -7890.058999999999
$ws.Range("N62").Value = -2748
$ws.Range("H65").Value = 8124.3887
$ws.Range("I65").Value = 8514.058999999999
$ws.Range("J65").Value = 1500
$ws.Range("K65").Value = 42570.295
$ws.Range("L65").Value = 7500
$ws.Range("M65").Value = -39450.295
$ws.Range("N65").Value = -13740
$ws.Range("H99").Value = 5895.6665
$ws.Range("I99").Value = 5293.75
$ws.Range("J99").Value = 7099.5
$ws.Range("K99").Value = 15881.25
$ws.Range("L99").Value = 21298.5
$ws.Range("M99").Value = -14383.25
$ws.Range("N99").Value = -24294.5
$ws.Range("H106").Value = 3855.1177
$ws.Range("J106").Value = 6661
$ws.Range("L106").Value = 6661
$ws.Range("N106").Value = -7923
$ws.Range("H127").Value = 2132.318
$ws.Range("J127").Value = 5601.6665
$ws.Range("L127").Value = 16804.9995
$ws.Range("N127").Value = -26724.9995
$ws.Range("H137").Value = 4626.8887
$ws.Range("I137").Value = 4187.364
$ws.Range("J137").Value = 6560.8
$ws.Range("K137").Value = 12562.092
$ws.Range("L137").Value = 19682.4
$ws.Range("M137").Value = -10012.092
$ws.Range("N137").Value = -24782.4
$ws.Range("H138").Value = 7901.9375
$ws.Range("J138").Value = 7860.338
$ws.Range("L138").Value = 23581.014
$ws.Range("N138").Value = -33861.014
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 371.76
$ws.Range("I5").Value = 414.2
$ws.Range("J5").Value = 202
$ws.Range("K5").Value = 414.2
$ws.Range("L5").Value = 202
$ws.Range("M5").Value = -302.2
$ws.Range("N5").Value = -426
$ws.Range("H32").Value = 8969.273999999999
$ws.Range("I32").Value = 8969.273999999999
$ws.Range("K32").Value = 8969.273999999999
$ws.Range("M32").Value = -8682.273999999999
$ws.Range("H74").Value = 5112.619
$ws.Range("I74").Value = 3903.8125
$ws.Range("K74").Value = 3903.8125
$ws.Range("M74").Value = -3029.8125
$ws.Range("H77").Value = 5112.619
$ws.Range("I77").Value = 3903.8125
$ws.Range("K77").Value = 19519.0625
$ws.Range("M77").Value = -15151.0625
$ws.Range("H97").Value = 1586.15
$ws.Range("I97").Value = 1662.2632
$ws.Range("K97").Value = 1662.2632
$ws.Range("M97").Value = -1166.2632
$ws.Range("H132").Value = 4471
$ws.Range("I132").Value = 3397.7368
$ws.Range("K132").Value = 10193.2104
$ws.Range("M132").Value = -7663.2104
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 371.76
$ws.Range("I4").Value = 414.2
$ws.Range("J4").Value = 202
$ws.Range("K4").Value = 414.2
$ws.Range("L4").Value = 202
$ws.Range("M4").Value = -299.2
$ws.Range("N4").Value = -432
$ws.Range("H20").Value = 4077.6667
$ws.Range("I20").Value = 3671.4285
$ws.Range("K20").Value = 3671.4285
$ws.Range("M20").Value = -3424.4285
$ws.Range("H94").Value = 1973.5
$ws.Range("I94").Value = 1164.8334
$ws.Range("J94").Value = 4399.5
$ws.Range("K94").Value = 1164.8334
$ws.Range("L94").Value = 4399.5
$ws.Range("M94").Value = -713.8334
$ws.Range("N94").Value = -5301.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1530.3077
$ws.Range("I22").Value = 466.33334
$ws.Range("J22").Value = 2442.2856
$ws.Range("K22").Value = 466.33334
$ws.Range("L22").Value = 2442.2856
$ws.Range("M22").Value = -116.33334
$ws.Range("N22").Value = -3142.2856
$ws.Range("H51").Value = 46999
$ws.Range("J51").Value = 46999
$ws.Range("L51").Value = 46999
$ws.Range("N51").Value = -48471
$ws.Range("H61").Value = 46999
$ws.Range("J61").Value = 46999
$ws.Range("L61").Value = 46999
$ws.Range("N61").Value = -47695
$ws.Range("H62").Value = 35717212
$ws.Range("J62").Value = 50003700
$ws.Range("L62").Value = 50003700
$ws.Range("N62").Value = -50004948
$ws.Range("H65").Value = 35717212
$ws.Range("J65").Value = 50003700
$ws.Range("L65").Value = 250018500
$ws.Range("N65").Value = -250024740
$ws.Range("H74").Value = 56467.89
$ws.Range("J74").Value = 58315.855
$ws.Range("L74").Value = 58315.855
$ws.Range("N74").Value = -60063.855
$ws.Range("H77").Value = 56467.89
$ws.Range("J77").Value = 58315.855
$ws.Range("L77").Value = 174947.565
$ws.Range("N77").Value = -183683.565
$ws.Range("H107").Value = 2220.8845
$ws.Range("I107").Value = 2578.875
$ws.Range("J107").Value = 1648.1
$ws.Range("K107").Value = 2578.875
$ws.Range("L107").Value = 1648.1
$ws.Range("M107").Value = -658.875
$ws.Range("N107").Value = -5488.1
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 232
$ws.Range("I23").Value = 825
$ws.Range("J23").Value = 157.875
$ws.Range("K23").Value = 2475
$ws.Range("L23").Value = 473.625
$ws.Range("M23").Value = -2240
$ws.Range("N23").Value = -943.625
$ws.Range("H34").Value = 145.25
$ws.Range("I34").Value = 145.25
$ws.Range("K34").Value = 435.75
$ws.Range("M34").Value = -351.75
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("N39").Value = 0
$ws.Range("L39").ClearContents()
$ws.Range("H55").Value = 4280
$ws.Range("J55").Value = 6666.6665
$ws.Range("L55").Value = 19999.9995
$ws.Range("N55").Value = -20353.9995
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3749.75
$ws.Range("I80").Value = 3499.5
$ws.Range("K80").Value = 3499.5
$ws.Range("M80").Value = -2501.5
$ws.Range("H83").Value = 3749.75
$ws.Range("I83").Value = 3499.5
$ws.Range("K83").Value = 17497.5
$ws.Range("M83").Value = -12505.5
$ws.Range("H97").Value = 433.33334
$ws.Range("I97").Value = 300
$ws.Range("J97").Value = 500
$ws.Range("K97").Value = 300
$ws.Range("L97").Value = 500
$ws.Range("N97").Value = -1492
$ws.Range("H126").Value = 12529.117
$ws.Range("I126").Value = 11999.833
$ws.Range("J126").Value = 13799.4
$ws.Range("K126").Value = 35999.499
$ws.Range("L126").Value = 41398.2
$ws.Range("M126").Value = -33529.499
$ws.Range("N126").Value = -46338.2
$ws.Range("H138").Value = 213994.5
$ws.Range("J138").Value = 213994.5
$ws.Range("L138").Value = 213994.5
$ws.Range("N138").Value = -224274.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5474.3794
$ws.Range("I61").Value = 2443.5
$ws.Range("K61").Value = 2443.5
$ws.Range("M61").Value = -2241.5
$ws.Range("H68").Value = 8197.777
$ws.Range("J68").Value = 10797.167
$ws.Range("L68").Value = 10797.167
$ws.Range("N68").Value = -12295.167
$ws.Range("H71").Value = 8197.777
$ws.Range("J71").Value = 10797.167
$ws.Range("L71").Value = 53985.835
$ws.Range("N71").Value = -61473.835
$ws.Range("H82").Value = 866.25
$ws.Range("I82").Value = 704.2857
$ws.Range("J82").Value = 2000
$ws.Range("K82").Value = 704.2857
$ws.Range("L82").Value = 2000
$ws.Range("M82").Value = -343.2857
$ws.Range("N82").Value = -2722
$ws.Range("H85").Value = 866.25
$ws.Range("I85").Value = 704.2857
$ws.Range("J85").Value = 2000
$ws.Range("K85").Value = 704.2857
$ws.Range("L85").Value = 2000
$ws.Range("M85").Value = 543.7143
$ws.Range("N85").Value = -4496
$ws.Range("H100").Value = 6601
$ws.Range("J100").Value = 6758.154
$ws.Range("L100").Value = 6758.154
$ws.Range("N100").Value = -7840.154
$ws.Range("H113").Value = 5474.3794
$ws.Range("I113").Value = 2443.5
$ws.Range("K113").Value = 2443.5
$ws.Range("M113").Value = -273.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 13677.143
$ws.Range("I62").Value = 14082.333
$ws.Range("K62").Value = 14082.333
$ws.Range("M62").Value = -13458.333
$ws.Range("H65").Value = 13677.143
$ws.Range("I65").Value = 14082.333
$ws.Range("K65").Value = 70411.66500000001
$ws.Range("M65").Value = -67291.66500000001
$ws.Range("H107").Value = 4543.515
$ws.Range("I107").Value = 3890.6365
$ws.Range("J107").Value = 5849.273
$ws.Range("K107").Value = 11671.9095
$ws.Range("L107").Value = 17547.819
$ws.Range("M107").Value = -9751.9095
$ws.Range("N107").Value = -21387.819
$ws.Range("H132").Value = 6439.6343
$ws.Range("I132").Value = 6192.5947
$ws.Range("J132").Value = 8724.75
$ws.Range("K132").Value = 18577.7841
$ws.Range("L132").Value = 26174.25
$ws.Range("M132").Value = -16047.7841
$ws.Range("N132").Value = -31234.25
$ws.Range("H136").Value = 2670.3208
$ws.Range("I136").Value = 1780.4286
$ws.Range("J136").Value = 6068.091
$ws.Range("K136").Value = 5341.2858
$ws.Range("L136").Value = 18204.273
$ws.Range("M136").Value = -2791.2858
$ws.Range("N136").Value = -23304.273
